$wb = $excel.ActiveWorkbook

# 1. Rename the shared string "Blok mapen new" -> "Blok Naling B"
#    This text lives in sheet "tambahBlokdanKamar", cell C2.
$wsTambah = $wb.Worksheets.Item("tambahBlokdanKamar")
$wsTambah.Range("C2").Value = "Blok Naling B"

# 2. Copy the block M1:Q5 (nomorKamar..lamaHuni columns) from tambahBlokdanKamar
#    into editBlokdanKamar starting at R1, then replicate row 5 into row 6 so the
#    new columns match the existing 6-row layout (row 6 mirrors row 5, just like
#    the existing B:Q columns already do on that sheet).
$wsEdit = $wb.Worksheets.Item("editBlokdanKamar")

$wsTambah.Range("M1:Q5").Copy()
$wsEdit.Range("R1").PasteSpecial(-4104)
$wsEdit.Range("R5:V5").Copy()
$wsEdit.Range("R6:V6").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# 3. Restore selections / view state for both sheets
$wsTambah.Activate()
$wsTambah.Application.ActiveWindow.ScrollRow = 1
$wsTambah.Range("D12").Select()

$wsEdit.Activate()
$wsEdit.Range("S33").Select()
